$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.558271169662476
$ws.Range("B1").Value = 1.849210858345032
$ws.Range("C1").Value = 1.971100091934204
$ws.Range("D1").Value = 2.332136392593384
$ws.Range("E1").Value = 3.206775903701782
